$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 4.33
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 2.1
$ws.Range("AH2").Value = 17
$ws.Range("AO2").Value = 7.5
$ws.Range("AZ2").Value = 101
$ws.Range("I3").Value = 3.9
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AN3").Value = 3.75
$ws.Range("AS3").Value = 251
$ws.Range("BA3").Value = 151
$ws.Range("V5").Value = 1.67
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2.85
$ws.Range("I6").Value = 2.45
$ws.Range("K6").Value = 1.9
$ws.Range("L6").Value = 3.1
$ws.Range("O6").Value = 1.52
$ws.Range("P6").Value = 2.22
$ws.Range("Q6").Value = 2.47
$ws.Range("S6").Value = 1.52
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = 2.07
$ws.Range("V6").Value = 1.6
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 13.5
$ws.Range("Y6").Value = 11.75
$ws.Range("AB6").Value = 55
$ws.Range("AC6").Value = 6
$ws.Range("AD6").Value = 5.8
$ws.Range("AE6").Value = 18.5
$ws.Range("AF6").Value = 120
$ws.Range("AH6").Value = 6
$ws.Range("AI6").Value = 10.5
$ws.Range("AJ6").Value = 10.25
$ws.Range("AK6").Value = 27
$ws.Range("AN6").Value = 4.65
$ws.Range("AO6").Value = 17.5
$ws.Range("AP6").Value = 30
$ws.Range("AQ6").Value = 90
$ws.Range("AS6").Value = 500
$ws.Range("AT6").Value = 2.18
$ws.Range("AU6").Value = 7.7
$ws.Range("AV6").Value = 90
$ws.Range("AW6").Value = 4.1
$ws.Range("AX6").Value = 13.5
$ws.Range("AZ6").Value = 65
